$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 5
$ws.Range("H4").Value = 5

# Row 6
$ws.Range("E6").Value = 28
$ws.Range("F6").Value = 16
$ws.Range("H6").Value = 16

# Row 7
$ws.Range("F7").Value = 6
$ws.Range("H7").Value = 6

# Row 9
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = 4

# Row 15
$ws.Range("E15").Value = 52
$ws.Range("F15").Value = 25
$ws.Range("H15").Value = 25

# Row 16
$ws.Range("E16").Value = 198

# Row 17
$ws.Range("E17").Value = 11

# Row 18
$ws.Range("E18").Value = 44
$ws.Range("F18").Value = 15
$ws.Range("H18").Value = 15
